# Daily attendance processing - rotate "Recorded By" (column G) entries.
# For every data row, move the first comma-separated name/email in column G
# to the end of the list (left-rotate by one element). Lists with a single
# entry are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ", "

    if ($parts.Length -gt 1) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
        $cell.Value2 = $rotated
    }
}
